{"js": "// Replace each division-problem text in the table with its new value.\n// Mapping is old \"a\u00f7b=\" -> new \"c\u00f7d=\" taken from the authoritative diff.\nconst replacements = [\n  [\"647\u00f78=\", \"722\u00f77=\"],\n  [\"624\u00f79=\", \"730\u00f75=\"],\n  [\"364\u00f72=\", \"503\u00f76=\"],\n  [\"198\u00f75=\", \"655\u00f73=\"],\n  [\"383\u00f75=\", \"556\u00f74=\"],\n  [\"918\u00f77=\", \"660\u00f79=\"],\n  [\"557\u00f78=\", \"770\u00f78=\"],\n  [\"285\u00f75=\", \"293\u00f77=\"],\n  [\"714\u00f73=\", \"330\u00f77=\"],\n  [\"571\u00f72=\", \"142\u00f72=\"],\n  [\"947\u00f76=\", \"618\u00f74=\"],\n  [\"603\u00f76=\", \"726\u00f77=\"],\n  [\"203\u00f77=\", \"749\u00f76=\"],\n  [\"971\u00f79=\", \"878\u00f74=\"],\n  [\"806\u00f79=\", \"947\u00f78=\"],\n  [\"216\u00f77=\", \"201\u00f74=\"],\n  [\"513\u00f76=\", \"568\u00f79=\"],\n  [\"651\u00f75=\", \"698\u00f72=\"],\n  [\"389\u00f73=\", \"996\u00f78=\"],\n  [\"693\u00f78=\", \"120\u00f73=\"],\n  [\"484\u00f73=\", \"144\u00f72=\"],\n  [\"495\u00f73=\", \"430\u00f76=\"],\n  [\"203\u00f74=\", \"890\u00f76=\"],\n  [\"697\u00f75=\", \"814\u00f77=\"],\n  [\"961\u00f72=\", \"123\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each division-problem text in the table with its new value.\n# Mapping is old \"a\u00f7b=\" -> new \"c\u00f7d=\" taken from the authoritative diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"647\u00f78=\", \"722\u00f77=\"),\n    @(\"624\u00f79=\", \"730\u00f75=\"),\n    @(\"364\u00f72=\", \"503\u00f76=\"),\n    @(\"198\u00f75=\", \"655\u00f73=\"),\n    @(\"383\u00f75=\", \"556\u00f74=\"),\n    @(\"918\u00f77=\", \"660\u00f79=\"),\n    @(\"557\u00f78=\", \"770\u00f78=\"),\n    @(\"285\u00f75=\", \"293\u00f77=\"),\n    @(\"714\u00f73=\", \"330\u00f77=\"),\n    @(\"571\u00f72=\", \"142\u00f72=\"),\n    @(\"947\u00f76=\", \"618\u00f74=\"),\n    @(\"603\u00f76=\", \"726\u00f77=\"),\n    @(\"203\u00f77=\", \"749\u00f76=\"),\n    @(\"971\u00f79=\", \"878\u00f74=\"),\n    @(\"806\u00f79=\", \"947\u00f78=\"),\n    @(\"216\u00f77=\", \"201\u00f74=\"),\n    @(\"513\u00f76=\", \"568\u00f79=\"),\n    @(\"651\u00f75=\", \"698\u00f72=\"),\n    @(\"389\u00f73=\", \"996\u00f78=\"),\n    @(\"693\u00f78=\", \"120\u00f73=\"),\n    @(\"484\u00f73=\", \"144\u00f72=\"),\n    @(\"495\u00f73=\", \"430\u00f76=\"),\n    @(\"203\u00f74=\", \"890\u00f76=\"),\n    @(\"697\u00f75=\", \"814\u00f77=\"),\n    @(\"961\u00f72=\", \"123\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $newText, 2)  # 2 = wdReplaceAll\n}\n"}
